# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

# row -> new value for column F
$updates = @{
    3  = 104
    4  = 1586
    5  = 604
    6  = 1090
    7  = 5
    8  = 11360
    9  = 19
    13 = 1087
    15 = 12331
    16 = 12990
    18 = 138
    23 = 93
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
